$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A before the existing Code/Description/Definition columns
$ws.Range("A1").EntireColumn.Insert()

# Treat the new column as text so "1.0" is stored as a string, not a number
$ws.Range("A1:A5").NumberFormat = "@"

# New header for the inserted column
$ws.Range("A1").Value = "Version"

# New values for rows 2-5 in the inserted column
$ws.Range("A2").Value = "1.0"
$ws.Range("A3").Value = "1.0"
$ws.Range("A4").Value = "1.0"
$ws.Range("A5").Value = "1.0"

# sheetFormatPr gains baseColWidth="10"
$ws.StandardWidth = 10
